$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prefix each category's sub-item labels with the category name ---

# Civilian
$ws.Cells.Item(7, 1).Value = "     Civilian, New nominations"
$ws.Cells.Item(8, 1).Value = "     Civilian, Carryover nominations"
$ws.Cells.Item(9, 1).Value = "     Civilian, Confirmed "
$ws.Cells.Item(10, 1).Value = "     Civilian, Unconfirmed "
$ws.Cells.Item(11, 1).Value = "     Civilian, Withdrawn "
$ws.Cells.Item(12, 1).Value = "     Civilian, Returned to White House "

# Other Civilian
$ws.Cells.Item(14, 1).Value = "     Other Civilian, New nominations"
$ws.Cells.Item(15, 1).Value = "     Other Civilian, Carryover nominations"
$ws.Cells.Item(16, 1).Value = "     Other Civilian, Confirmed "
$ws.Cells.Item(17, 1).Value = "     Other Civilian, Unconfirmed "

# Air Force
$ws.Cells.Item(19, 1).Value = "     Air Force, New nominations"
$ws.Cells.Item(20, 1).Value = "     Air Force, Carryover nominations"
$ws.Cells.Item(21, 1).Value = "     Air Force, Confirmed "
$ws.Cells.Item(22, 1).Value = "     Air Force, Returned to White House "

# Army
$ws.Cells.Item(24, 1).Value = "     Army, New nominations"
$ws.Cells.Item(25, 1).Value = "     Army, Carryover nominations"
$ws.Cells.Item(26, 1).Value = "     Army, Confirmed "
$ws.Cells.Item(27, 1).Value = "     Army, Unconfirmed "
$ws.Cells.Item(28, 1).Value = "     Army, Returned to White House "

# Navy
$ws.Cells.Item(30, 1).Value = "     Navy, New nominations"
$ws.Cells.Item(31, 1).Value = "     Navy, Carryover nominations"
$ws.Cells.Item(32, 1).Value = "     Navy, Confirmed "
$ws.Cells.Item(33, 1).Value = "     Navy, Returned to White House "

# Marine Corps
$ws.Cells.Item(35, 1).Value = "     Marine Corps, New nominations"
$ws.Cells.Item(36, 1).Value = "     Marine Corps, Carryover nominations"
$ws.Cells.Item(37, 1).Value = "     Marine Corps, Confirmed "

# --- Rework the Summary block (rows 38-44) into a new 6-row totals block (rows 38-43) ---

# Row 38 was the "Summary" header with no value; it becomes "Total new nominations"
# with the value that used to live in row 40 ("Total nominations received this session").
# Borrow the number style (s=3, #,##0 right aligned) from row 39's B cell, which already
# carries that style, so no new style entries get created.
$ws.Cells.Item(39, 2).Copy() | Out-Null
$ws.Cells.Item(38, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(38, 1).Value = "Total new nominations"
$ws.Cells.Item(38, 2).Value = 22162

# Row 39: "Total nominations carried over from First Session " -> "Total carryover nominations"
# (value 1150 is unchanged)
$ws.Cells.Item(39, 1).Value = "Total carryover nominations"

# Row 40: "Total nominations received this session " -> "Total confirmed " (value becomes 22512)
$ws.Cells.Item(40, 1).Value = "Total confirmed "
$ws.Cells.Item(40, 2).Value = 22512

# Row 41: "Total confirmed " -> "Total unconfirmed " (value becomes 21).
# This row's number style also changes from s=3 to s=2 (plain number, no #,##0 format),
# so borrow formatting from row 42's B cell, which already has that plain style.
$ws.Cells.Item(42, 2).Copy() | Out-Null
$ws.Cells.Item(41, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(41, 1).Value = "Total unconfirmed "
$ws.Cells.Item(41, 2).Value = 21

# Row 42: "Total unconfirmed " -> "Total withdrawn " (value becomes 13)
$ws.Cells.Item(42, 1).Value = "Total withdrawn "
$ws.Cells.Item(42, 2).Value = 13

# Row 43: "Total withdrawn " -> "Total returned" (value becomes 766, taken from the old row 44)
$ws.Cells.Item(43, 1).Value = "Total returned"
$ws.Cells.Item(43, 2).Value = 766

# Old row 44 ("Total Returned to White House " / 766) is no longer needed; remove it entirely.
$ws.Rows.Item(44).Delete()
